# Generate Report for Archive
# 1) Update the localization status text everywhere it appears
#    ("Ready for handoff" -> "In Translation") across all sheets.
# 2) Shrink the "Status"/language columns (Overview!E:F and the
#    per-language "Status" column C on the zh-cn / de-de sheets) from
#    17.2159881591797 chars wide down to 13.4101845877511 chars wide.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$xlWhole = 1

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldStatus, $newStatus, $xlWhole) | Out-Null
}

# Column width change.
# The ColumnWidth COM property is expressed in characters and Excel
# snaps the stored width to the nearest 1/6th of a character (pixel
# grid), so feed it the value whose rounded result lands closest to
# the target stored width of 13.4101845877511.
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
